# Generate Report for Handoff
# Update status + timestamps to reflect a fresh handoff generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn sheet: Status (C2) + Latest Handoff Datetime (H2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-06 20:53:34"

# de-de sheet: Status (C2) + Latest Handoff Datetime (H2)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-06 20:53:39"

# Overview sheet: zh-cn status (E2), de-de status (F2), Latest HO Xliff Generate Date (G2)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-06 20:53:39"

# Column widths widen slightly on the Status columns after the text change (auto-fit
# to fit the longer "Ready for handoff" string).
$wsOverview.Columns("E:F").ColumnWidth = 16.33
$wsZhCn.Columns("C:C").ColumnWidth = 16.33
$wsDeDe.Columns("C:C").ColumnWidth = 16.33
